$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column Z with header "13-jul" and matching daily values
$ws.Range("Z1").Value = "13-jul"
$ws.Range("Z1").NumberFormat = "@"

$ws.Range("Z2").Value = 13
$ws.Range("Z2").HorizontalAlignment = -4108
$ws.Range("Z2").NumberFormat = "0"

$ws.Range("Z3").Value = 20
$ws.Range("Z3").HorizontalAlignment = -4108
$ws.Range("Z3").NumberFormat = "0"

$ws.Range("Z4").Value = 8
$ws.Range("Z4").HorizontalAlignment = -4108
$ws.Range("Z4").NumberFormat = "0"

$ws.Range("Z5").Value = 8
$ws.Range("Z5").HorizontalAlignment = -4108
$ws.Range("Z5").NumberFormat = "0"

$ws.Range("Z6").Value = 12
$ws.Range("Z6").HorizontalAlignment = -4108
$ws.Range("Z6").NumberFormat = "0"

$ws.Range("Z7").Value = 16
$ws.Range("Z7").HorizontalAlignment = -4108
$ws.Range("Z7").NumberFormat = "0"

$ws.Range("Z8").Value = 18
$ws.Range("Z8").HorizontalAlignment = -4108
$ws.Range("Z8").NumberFormat = "0"

$ws.Range("Z9").Value = 12
$ws.Range("Z9").HorizontalAlignment = -4108
$ws.Range("Z9").NumberFormat = "0"

$ws.Range("Z10").Value = 19
$ws.Range("Z10").HorizontalAlignment = -4108
$ws.Range("Z10").NumberFormat = "0"

$ws.Range("Z11").Value = 24
$ws.Range("Z11").HorizontalAlignment = -4108
$ws.Range("Z11").NumberFormat = "0"

# Update selection to match target state
$ws.Range("V11").Select()
